# Add a new "30-ago" column (BI) to the worksheet, mirroring the last
# existing "29-ago" column (BH), and fill in new values for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same style as BH1 (text custom format used for the other
# date-label headers), new label "30-ago".
$ws.Range("BI1").Value = "30-ago"
$ws.Range("BI1").NumberFormat = $ws.Range("BH1").NumberFormat

# New numeric values for BI2:BI18, one per data row.
$values = @{
    2  = 0
    3  = 11.790335228036453
    4  = 13.289322258264672
    5  = 13.653018835169508
    6  = 0
    7  = 3.9227148084470347
    8  = 10.508508770588659
    9  = 5.2889086594976442
    10 = 17.840837908014674
    11 = 13.387885519631487
    12 = 0
    13 = 9.1914140751933608
    14 = 0
    15 = 0
    16 = 18.926216224121152
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 61).Value = $values[$row]
}

# Mirror the selection change recorded in the saved file.
$ws.Range("BJ4").Select()
